$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# Add the new "metadata" sheet right after "data"
$ws = $wb.Worksheets.Add($null, $dataSheet)
$ws.Name = "metadata"

# Copy the header-row cell formatting (bold, centered, bordered) from data!B1
# onto the metadata header row B1:G1 so it reuses the same style record.
$dataSheet.Range("B1").Copy()
$ws.Range("B1:G1").PasteSpecial(-4122)

# Copy the index-column cell formatting from data!A2 onto metadata!A2.
$dataSheet.Range("A2").Copy()
$ws.Range("A2").PasteSpecial(-4122)

$ws.Application.CutCopyMode = $false

# Header row
$ws.Range("B1").Value = "data_name"
$ws.Range("C1").Value = "data_id"
$ws.Range("D1").Value = "data_version"
$ws.Range("E1").Value = "data_version_created"
$ws.Range("F1").Value = "panel_query_time"
$ws.Range("G1").Value = "panel_get_request"

# Data row
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "Malignant Hyperthermia Susceptibility"
$ws.Range("C2").Value = 3378

# data_version ("1.3") must land as TEXT, not a number - build it via a
# formula producing a string, then paste back as a value so no stray style
# (e.g. quotePrefix) is left behind on the cell.
$ws.Range("D2").Formula = "=""1.3"""
$ws.Range("D2").Copy()
$ws.Range("D2").PasteSpecial(-4163)
$ws.Application.CutCopyMode = $false

$ws.Range("E2").Value = "2021-01-29T02:54:23.037896Z"
$ws.Range("F2").Value = "2021-10-05 14:34:33.423892"
$ws.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/3378/?format=json"

# Refresh timestamps on the "data" sheet (F2:F7) per the updated panel query time
$dataSheet.Range("F2").Value = "2021-10-05 14:34:33.427820"
$dataSheet.Range("F3").Value = "2021-10-05 14:34:33.427828"
$dataSheet.Range("F4").Value = "2021-10-05 14:34:33.427831"
$dataSheet.Range("F5").Value = "2021-10-05 14:34:33.427833"
$dataSheet.Range("F6").Value = "2021-10-05 14:34:33.427836"
$dataSheet.Range("F7").Value = "2021-10-05 14:34:33.427839"

[void]$dataSheet.Activate()
